# ------------------------------------------------------------------
# Edit: add a "metadata" worksheet (with panel export metadata) after
# the existing "data" sheet, and refresh the per-row "time_taken"
# timestamps on the "data" sheet to reflect the new export run.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh "time_taken" (column F) timestamps on the data sheet ---
$timeTaken = @(
    "2021-10-05 14:33:16.758401",
    "2021-10-05 14:33:16.758409",
    "2021-10-05 14:33:16.758412",
    "2021-10-05 14:33:16.758415",
    "2021-10-05 14:33:16.758418",
    "2021-10-05 14:33:16.758421",
    "2021-10-05 14:33:16.758423",
    "2021-10-05 14:33:16.758426",
    "2021-10-05 14:33:16.758429",
    "2021-10-05 14:33:16.758432",
    "2021-10-05 14:33:16.758434",
    "2021-10-05 14:33:16.758437",
    "2021-10-05 14:33:16.758439",
    "2021-10-05 14:33:16.758441",
    "2021-10-05 14:33:16.758444",
    "2021-10-05 14:33:16.758446",
    "2021-10-05 14:33:16.758449",
    "2021-10-05 14:33:16.758452",
    "2021-10-05 14:33:16.758455",
    "2021-10-05 14:33:16.758457",
    "2021-10-05 14:33:16.758460",
    "2021-10-05 14:33:16.758463",
    "2021-10-05 14:33:16.758465",
    "2021-10-05 14:33:16.758468",
    "2021-10-05 14:33:16.758471",
    "2021-10-05 14:33:16.758474",
    "2021-10-05 14:33:16.758477",
    "2021-10-05 14:33:16.758479",
    "2021-10-05 14:33:16.758482",
    "2021-10-05 14:33:16.758484",
    "2021-10-05 14:33:16.758487",
    "2021-10-05 14:33:16.758490",
    "2021-10-05 14:33:16.758493",
    "2021-10-05 14:33:16.758496",
    "2021-10-05 14:33:16.758498",
    "2021-10-05 14:33:16.758501",
    "2021-10-05 14:33:16.758503",
    "2021-10-05 14:33:16.758506",
    "2021-10-05 14:33:16.758509",
    "2021-10-05 14:33:16.758511",
    "2021-10-05 14:33:16.758514",
    "2021-10-05 14:33:16.758517",
    "2021-10-05 14:33:16.758519",
    "2021-10-05 14:33:16.758522",
    "2021-10-05 14:33:16.758525",
    "2021-10-05 14:33:16.758527",
    "2021-10-05 14:33:16.758530",
    "2021-10-05 14:33:16.758532",
    "2021-10-05 14:33:16.758535",
    "2021-10-05 14:33:16.758538",
    "2021-10-05 14:33:16.758540",
    "2021-10-05 14:33:16.758543",
    "2021-10-05 14:33:16.758546",
    "2021-10-05 14:33:16.758549",
    "2021-10-05 14:33:16.758552",
    "2021-10-05 14:33:16.758554",
    "2021-10-05 14:33:16.758557",
    "2021-10-05 14:33:16.758560",
    "2021-10-05 14:33:16.758562",
    "2021-10-05 14:33:16.758565",
    "2021-10-05 14:33:16.758568",
    "2021-10-05 14:33:16.758570",
    "2021-10-05 14:33:16.758573",
    "2021-10-05 14:33:16.758576",
    "2021-10-05 14:33:16.758580",
    "2021-10-05 14:33:16.758583",
    "2021-10-05 14:33:16.758585",
    "2021-10-05 14:33:16.758588",
    "2021-10-05 14:33:16.758591",
    "2021-10-05 14:33:16.758593",
    "2021-10-05 14:33:16.758596",
    "2021-10-05 14:33:16.758599",
    "2021-10-05 14:33:16.758601",
    "2021-10-05 14:33:16.758604",
    "2021-10-05 14:33:16.758606",
    "2021-10-05 14:33:16.758609",
    "2021-10-05 14:33:16.758613",
    "2021-10-05 14:33:16.758617",
    "2021-10-05 14:33:16.758619",
    "2021-10-05 14:33:16.758622",
    "2021-10-05 14:33:16.758625",
    "2021-10-05 14:33:16.758628",
    "2021-10-05 14:33:16.758630",
    "2021-10-05 14:33:16.758633",
    "2021-10-05 14:33:16.758636",
    "2021-10-05 14:33:16.758638",
    "2021-10-05 14:33:16.758641",
    "2021-10-05 14:33:16.758643",
    "2021-10-05 14:33:16.758646",
    "2021-10-05 14:33:16.758649",
    "2021-10-05 14:33:16.758651",
    "2021-10-05 14:33:16.758654",
    "2021-10-05 14:33:16.758658",
    "2021-10-05 14:33:16.758661",
    "2021-10-05 14:33:16.758663",
    "2021-10-05 14:33:16.758666",
    "2021-10-05 14:33:16.758668",
    "2021-10-05 14:33:16.758671",
    "2021-10-05 14:33:16.758674",
    "2021-10-05 14:33:16.758676",
    "2021-10-05 14:33:16.758679",
    "2021-10-05 14:33:16.758681",
    "2021-10-05 14:33:16.758684",
    "2021-10-05 14:33:16.758687",
    "2021-10-05 14:33:16.758689",
    "2021-10-05 14:33:16.758692",
    "2021-10-05 14:33:16.758695",
    "2021-10-05 14:33:16.758698",
    "2021-10-05 14:33:16.758702",
    "2021-10-05 14:33:16.758706",
    "2021-10-05 14:33:16.758708",
    "2021-10-05 14:33:16.758711",
    "2021-10-05 14:33:16.758714",
    "2021-10-05 14:33:16.758716",
    "2021-10-05 14:33:16.758719",
    "2021-10-05 14:33:16.758721",
    "2021-10-05 14:33:16.758724",
    "2021-10-05 14:33:16.758727",
    "2021-10-05 14:33:16.758729"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# --- 2. Add the new "metadata" sheet, placed right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the outline properties openpyxl stamps on every sheet ($sheetPr)
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Reuse the bold/border/center header style (style index 1) already
# defined on the "data" sheet header row, instead of re-deriving a new
# style via Font/Border/Alignment property writes.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# --- Header row ---
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# --- Data row ---
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Bleeding and Platelet Disorders"
$metaSheet.Range("C2").Value = 54

# Force "1.2" to be stored as text (not the number 1.2): write it as a
# formula producing a text result, then collapse the formula down to a
# static value-only paste so no extra number-format style gets created.
$metaSheet.Range("D2").Formula = '="1.2"'
$metaSheet.Range("D2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)

$metaSheet.Range("E2").Value = "2021-06-15T10:04:16.914690Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:16.755016"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/54/?format=json"

# Keep "data" the active/selected sheet (matches original activeTab=0)
$dataSheet.Activate()
$dataSheet.Range("A1").Select()

Write-Host "edit applied"
